# Append two new trading-log rows (58 and 59) to Sheet1, mirroring the
# existing TRADING_ATTEMPT / POSITION_FAILED pair pattern already present
# in the log (e.g. rows 56-57) for a new SUI attempt.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58: TRADING_ATTEMPT
$ws.Cells.Item(58, 1).Value = "2025-10-03T12:41:58.629075"
$ws.Cells.Item(58, 2).Value = "TRADING_ATTEMPT"
$ws.Cells.Item(58, 3).Value = "SUI"
$ws.Cells.Item(58, 4).Value = "UNKNOWN"
$ws.Cells.Item(58, 5).Value = 3.582294139231992
$ws.Cells.Item(58, 11).Value = "ATTEMPT"
$ws.Cells.Item(58, 12).Value = "Attempting trade 1/1"

# Row 59: POSITION_FAILED
$ws.Cells.Item(59, 1).Value = "2025-10-03T12:42:00.542078"
$ws.Cells.Item(59, 2).Value = "POSITION_FAILED"
$ws.Cells.Item(59, 3).Value = "SUI"
$ws.Cells.Item(59, 4).Value = "UNKNOWN"
$ws.Cells.Item(59, 11).Value = "FAILED"
$ws.Cells.Item(59, 12).Value = "Trade execution failed for trade 1"
